# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor TPM-derived statistics:
# refresh expression/specificity values (columns E-J, M-T) across rows 2-17
# to reflect the updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 1.782436333333333
$ws.Cells.Item(2, 8).Value = 5.347308999999999
$ws.Cells.Item(2, 9).Value = 0.00914036392049929
$ws.Cells.Item(2, 10).Value = 0.009140363920499292
$ws.Cells.Item(2, 13).Value = 145.7007446666667
$ws.Cells.Item(2, 14).Value = 437.1022340000001
$ws.Cells.Item(2, 15).Value = 0.2865937750105843
$ws.Cells.Item(2, 16).Value = 0.2865937750105843
$ws.Cells.Item(2, 17).Value = 259.7023010875895
$ws.Cells.Item(2, 18).Value = 2337.320709788306
$ws.Cells.Item(2, 19).Value = 0.002619571400946435
$ws.Cells.Item(2, 20).Value = 0.002619571400946436

$ws.Cells.Item(3, 7).Value = 1.782436333333333
$ws.Cells.Item(3, 8).Value = 5.347308999999999
$ws.Cells.Item(3, 9).Value = 0.00914036392049929
$ws.Cells.Item(3, 10).Value = 0.009140363920499292
$ws.Cells.Item(3, 15).Value = 0.3320294904365841
$ws.Cells.Item(3, 16).Value = 0.3320294904365841
$ws.Cells.Item(3, 17).Value = 300.8747230889302
$ws.Cells.Item(3, 18).Value = 2707.872507800372
$ws.Cells.Item(3, 19).Value = 0.003034870374928318
$ws.Cells.Item(3, 20).Value = 0.003034870374928318

$ws.Cells.Item(4, 7).Value = 1.782436333333333
$ws.Cells.Item(4, 8).Value = 5.347308999999999
$ws.Cells.Item(4, 9).Value = 0.00914036392049929
$ws.Cells.Item(4, 10).Value = 0.009140363920499292
$ws.Cells.Item(4, 13).Value = 128.1261546666667
$ws.Cells.Item(4, 14).Value = 384.378464
$ws.Cells.Item(4, 15).Value = 0.2520245069956105
$ws.Cells.Item(4, 16).Value = 0.2520245069956105
$ws.Cells.Item(4, 17).Value = 228.3767133281529
$ws.Cells.Item(4, 18).Value = 2055.390419953376
$ws.Cells.Item(4, 19).Value = 0.002303595710824299
$ws.Cells.Item(4, 20).Value = 0.002303595710824299

$ws.Cells.Item(5, 7).Value = 1.782436333333333
$ws.Cells.Item(5, 8).Value = 5.347308999999999
$ws.Cells.Item(5, 9).Value = 0.00914036392049929
$ws.Cells.Item(5, 10).Value = 0.009140363920499292
$ws.Cells.Item(5, 13).Value = 65.761079
$ws.Cells.Item(5, 14).Value = 197.283237
$ws.Cells.Item(5, 15).Value = 0.1293522275572212
$ws.Cells.Item(5, 16).Value = 0.1293522275572212
$ws.Cells.Item(5, 17).Value = 117.2149365288036
$ws.Cells.Item(5, 18).Value = 1054.934428759233
$ws.Cells.Item(5, 19).Value = 0.001182326433800239
$ws.Cells.Item(5, 20).Value = 0.001182326433800239

$ws.Cells.Item(6, 8).Value = 564.692825
$ws.Cells.Item(6, 9).Value = 0.965251479537618
$ws.Cells.Item(6, 10).Value = 0.965251479537618
$ws.Cells.Item(6, 13).Value = 145.7007446666667
$ws.Cells.Item(6, 14).Value = 437.1022340000001
$ws.Cells.Item(6, 15).Value = 0.2865937750105843
$ws.Cells.Item(6, 16).Value = 0.2865937750105843
$ws.Cells.Item(6, 17).Value = 27425.38837014123
$ws.Cells.Item(6, 18).Value = 246828.4953312711
$ws.Cells.Item(6, 19).Value = 0.2766350653552377
$ws.Cells.Item(6, 20).Value = 0.2766350653552377

$ws.Cells.Item(7, 8).Value = 564.692825
$ws.Cells.Item(7, 9).Value = 0.965251479537618
$ws.Cells.Item(7, 10).Value = 0.965251479537618
$ws.Cells.Item(7, 15).Value = 0.3320294904365841
$ws.Cells.Item(7, 16).Value = 0.3320294904365841
$ws.Cells.Item(7, 17).Value = 31773.32698600001
$ws.Cells.Item(7, 18).Value = 285959.9428740001
$ws.Cells.Item(7, 19).Value = 0.3204919568940343
$ws.Cells.Item(7, 20).Value = 0.3204919568940343

$ws.Cells.Item(8, 8).Value = 564.692825
$ws.Cells.Item(8, 9).Value = 0.965251479537618
$ws.Cells.Item(8, 10).Value = 0.965251479537618
$ws.Cells.Item(8, 13).Value = 128.1261546666667
$ws.Cells.Item(8, 14).Value = 384.378464
$ws.Cells.Item(8, 15).Value = 0.2520245069956105
$ws.Cells.Item(8, 16).Value = 0.2520245069956105
$ws.Cells.Item(8, 17).Value = 24117.30674503565
$ws.Cells.Item(8, 18).Value = 217055.7607053208
$ws.Cells.Item(8, 19).Value = 0.2432670282572518
$ws.Cells.Item(8, 20).Value = 0.2432670282572518

$ws.Cells.Item(9, 8).Value = 564.692825
$ws.Cells.Item(9, 9).Value = 0.965251479537618
$ws.Cells.Item(9, 10).Value = 0.965251479537618
$ws.Cells.Item(9, 13).Value = 65.761079
$ws.Cells.Item(9, 14).Value = 197.283237
$ws.Cells.Item(9, 15).Value = 0.1293522275572212
$ws.Cells.Item(9, 16).Value = 0.1293522275572212
$ws.Cells.Item(9, 17).Value = 12378.26982518606
$ws.Cells.Item(9, 18).Value = 111404.4284266745
$ws.Cells.Item(9, 19).Value = 0.1248574290310944
$ws.Cells.Item(9, 20).Value = 0.1248574290310944

$ws.Cells.Item(10, 7).Value = 4.870778333333333
$ws.Cells.Item(10, 8).Value = 14.612335
$ws.Cells.Item(10, 9).Value = 0.02497743437460768
$ws.Cells.Item(10, 10).Value = 0.02497743437460768
$ws.Cells.Item(10, 13).Value = 145.7007446666667
$ws.Cells.Item(10, 14).Value = 437.1022340000001
$ws.Cells.Item(10, 15).Value = 0.2865937750105843
$ws.Cells.Item(10, 16).Value = 0.2865937750105843
$ws.Cells.Item(10, 17).Value = 709.6760302729323
$ws.Cells.Item(10, 18).Value = 6387.08427245639
$ws.Cells.Item(10, 19).Value = 0.007158377207497946
$ws.Cells.Item(10, 20).Value = 0.007158377207497946

$ws.Cells.Item(11, 7).Value = 4.870778333333333
$ws.Cells.Item(11, 8).Value = 14.612335
$ws.Cells.Item(11, 9).Value = 0.02497743437460768
$ws.Cells.Item(11, 10).Value = 0.02497743437460768
$ws.Cells.Item(11, 15).Value = 0.3320294904365841
$ws.Cells.Item(11, 16).Value = 0.3320294904365841
$ws.Cells.Item(11, 17).Value = 822.185934421909
$ws.Cells.Item(11, 18).Value = 7399.67340979718
$ws.Cells.Item(11, 19).Value = 0.008293244807814208
$ws.Cells.Item(11, 20).Value = 0.008293244807814208

$ws.Cells.Item(12, 7).Value = 4.870778333333333
$ws.Cells.Item(12, 8).Value = 14.612335
$ws.Cells.Item(12, 9).Value = 0.02497743437460768
$ws.Cells.Item(12, 10).Value = 0.02497743437460768
$ws.Cells.Item(12, 13).Value = 128.1261546666667
$ws.Cells.Item(12, 14).Value = 384.378464
$ws.Cells.Item(12, 15).Value = 0.2520245069956105
$ws.Cells.Item(12, 16).Value = 0.2520245069956105
$ws.Cells.Item(12, 17).Value = 624.0740980837156
$ws.Cells.Item(12, 18).Value = 5616.666882753439
$ws.Cells.Item(12, 19).Value = 0.006294925584275714
$ws.Cells.Item(12, 20).Value = 0.006294925584275714

$ws.Cells.Item(13, 7).Value = 4.870778333333333
$ws.Cells.Item(13, 8).Value = 14.612335
$ws.Cells.Item(13, 9).Value = 0.02497743437460768
$ws.Cells.Item(13, 10).Value = 0.02497743437460768
$ws.Cells.Item(13, 13).Value = 65.761079
$ws.Cells.Item(13, 14).Value = 197.283237
$ws.Cells.Item(13, 15).Value = 0.1293522275572212
$ws.Cells.Item(13, 16).Value = 0.1293522275572212
$ws.Cells.Item(13, 17).Value = 320.3076387698216
$ws.Cells.Item(13, 18).Value = 2882.768748928394
$ws.Cells.Item(13, 19).Value = 0.003230886775019812
$ws.Cells.Item(13, 20).Value = 0.003230886775019812

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.1229953333333333
$ws.Cells.Item(14, 8).Value = 0.368986
$ws.Cells.Item(14, 9).Value = 0.0006307221672750447
$ws.Cells.Item(14, 10).Value = 0.0006307221672750447
$ws.Cells.Item(14, 13).Value = 145.7007446666667
$ws.Cells.Item(14, 14).Value = 437.1022340000001
$ws.Cells.Item(14, 15).Value = 0.2865937750105843
$ws.Cells.Item(14, 16).Value = 0.2865937750105843
$ws.Cells.Item(14, 17).Value = 17.92051165719156
$ws.Cells.Item(14, 18).Value = 161.284604914724
$ws.Cells.Item(14, 19).Value = 0.0001807610469022122
$ws.Cells.Item(14, 20).Value = 0.0001807610469022122

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.1229953333333333
$ws.Cells.Item(15, 8).Value = 0.368986
$ws.Cells.Item(15, 9).Value = 0.0006307221672750447
$ws.Cells.Item(15, 10).Value = 0.0006307221672750447
$ws.Cells.Item(15, 15).Value = 0.3320294904365841
$ws.Cells.Item(15, 16).Value = 0.3320294904365841
$ws.Cells.Item(15, 17).Value = 20.76157569605423
$ws.Cells.Item(15, 18).Value = 186.854181264488
$ws.Cells.Item(15, 19).Value = 0.0002094183598073911
$ws.Cells.Item(15, 20).Value = 0.0002094183598073911

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.1229953333333333
$ws.Cells.Item(16, 8).Value = 0.368986
$ws.Cells.Item(16, 9).Value = 0.0006307221672750447
$ws.Cells.Item(16, 10).Value = 0.0006307221672750447
$ws.Cells.Item(16, 13).Value = 128.1261546666667
$ws.Cells.Item(16, 14).Value = 384.378464
$ws.Cells.Item(16, 15).Value = 0.2520245069956105
$ws.Cells.Item(16, 16).Value = 0.2520245069956105
$ws.Cells.Item(16, 17).Value = 15.75891910194489
$ws.Cells.Item(16, 18).Value = 141.830271917504
$ws.Cells.Item(16, 19).Value = 0.0001589574432586961
$ws.Cells.Item(16, 20).Value = 0.0001589574432586961

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.1229953333333333
$ws.Cells.Item(17, 8).Value = 0.368986
$ws.Cells.Item(17, 9).Value = 0.0006307221672750447
$ws.Cells.Item(17, 10).Value = 0.0006307221672750447
$ws.Cells.Item(17, 13).Value = 65.761079
$ws.Cells.Item(17, 14).Value = 197.283237
$ws.Cells.Item(17, 15).Value = 0.1293522275572212
$ws.Cells.Item(17, 16).Value = 0.1293522275572212
$ws.Cells.Item(17, 17).Value = 8.088305831964666
$ws.Cells.Item(17, 18).Value = 72.794752487682
$ws.Cells.Item(17, 19).Value = 0.00008158531730674534
$ws.Cells.Item(17, 20).Value = 0.00008158531730674534
